$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.426.52'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '3.829.65'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''600.26'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = '''163.29'
$ws.Range("E6").Value = '  -2.85%  '
$ws.Range("D7").Value = '3.827.25'
$ws.Range("E7").Value = '  +2.08%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -2.36%  '
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").Value = '''6.32'
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("D12").Value = '''0.459'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '''36.78'
$ws.Range("E13").Value = '  -3.97%  '
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '4.469.05'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '3.829.85'
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("D17").Value = '68.617.65'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '''7.57'
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '''17.10'
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").Value = '''11.18'
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("D22").Value = '''485.40'
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("E24").Value = '  +6.60%  '
$ws.Range("D25").Value = '''84.05'
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("D27").Value = '''12.09'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").Value = '''9.99'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").Value = '3.982.16'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("E33").Value = '  -4.40%  '
$ws.Range("D34").Value = '''31.78'
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").Value = '3.774.61'
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("E37").Value = '  +1.12%  '
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("E39").Value = '  -1.93%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("D43").Value = '''428.45'
$ws.Range("E43").Value = '  +0.94%  '
$ws.Range("D44").Value = '''48.45'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("D45").Value = '''1.98'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D47").Value = '''8.41'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("D48").Value = '2.842.48'
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("D49").Value = '''142.65'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0357'
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''25.90'
$ws.Range("E51").Value = '  +12.54%  '
